# Implementação parcial do cadastro de patente
# Adds the new apropriação entry (row 14) to the sheet and updates the
# current selection to match, mirroring what a user would do after typing
# the new date/hours pair into the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
$ws.Range("A14").Value = 41567
$ws.Range("B14").Value = 0.1423611111111111

# Match the formatting already used by the rows above (date / duration).
[void]$ws.Range("A13").Copy()
[void]$ws.Range("A14").PasteSpecial(-4122)

[void]$ws.Range("B13").Copy()
[void]$ws.Range("B14").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Selection ----------------------------------------------------------
[void]$ws.Range("B4:B14").Select()
